$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header-like numeric row, columns B..J values 0..8 (style already s=1 on B1, same style applied to new cells)
$row1 = @(0,1,2,3,4,5,6,7,8)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $col = 2 + $i   # B=2 .. J=10
    $ws.Cells.Item(1, $col).Value = $row1[$i]
}

# Row 2: carID
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "carID"
$row2 = @(1,2,3,4,6,9,12,12)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = 3 + $i   # C=3 .. J=10
    $ws.Cells.Item(2, $col).Value = $row2[$i]
}

# Row 3: speed2
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "speed2"
$row3 = @(42.43,55.9,49.24,43.01,40.31,47.17,40.31,44.72)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}

# Row 4: asma
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "asma"
$row4 = @(21.22,59.72,40.7,22.89,15.18,34.77,15.18,27.78)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(4, $col).Value = $row4[$i]
}

# Row 5: ceza_tutar
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "ceza_tutar"
$row5 = @(1508.5,6440,3136,1508.5,1508.5,3136,1508.5,1508.5)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(5, $col).Value = $row5[$i]
}

# Row 6: hesaplanan_asma
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "hesaplanan_asma"
$row6 = @(10,50,30,10,10,30,10,10)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(6, $col).Value = $row6[$i]
}
